$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.240.79'
$ws.Range('E2').Value = '  +4.17%  '
$ws.Range('D3').Value = '3.637.17'
$ws.Range('E3').Value = '  +3.30%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Formula = '="203.31"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +12.15%  '
$ws.Range('D6').Formula = '="565.29"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').Value = '3.614.38'
$ws.Range('E7').Value = '  +2.94%  '
$ws.Range('E8').Value = '  +2.59%  '
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('E10').Value = '  +2.74%  '
$ws.Range('D11').Formula = '="57.83"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +7.79%  '
$ws.Range('D12').Formula = '="0.153"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  +8.34%  '
$ws.Range('D13').Formula = '="0.0000290"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +17.46%  '
$ws.Range('D14').Formula = '="10.02"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  +3.89%  '
$ws.Range('D15').Value = '4.213.14'
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '3.636.31'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '68.137.37'
$ws.Range('E18').Value = '  +4.40%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Formula = '="12.45"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +3.79%  '
$ws.Range('D20').Formula = '="18.54"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Formula = '="1.08"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').Formula = '="400.68"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +3.48%  '
$ws.Range('D23').Formula = '="12.97"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  +27.61%  '
$ws.Range('D24').Formula = '="4.17"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').Formula = '="85.59"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +2.61%  '
$ws.Range('E26').Value = '  +4.28%  '
$ws.Range('D27').Formula = '="12.56"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +3.02%  '
$ws.Range('E28').Value = '  +2.02%  '
$ws.Range('D29').Formula = '="3.83"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +8.37%  '
$ws.Range('D30').Formula = '="8.16"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +22.30%  '
$ws.Range('E31').Value = '  +4.01%  '
$ws.Range('D32').Formula = '="31.80"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +3.96%  '
$ws.Range('D33').Formula = '="695.52"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +14.66%  '
$ws.Range('D34').Formula = '="12.22"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +2.77%  '
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('D36').Formula = '="64.21"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').Formula = '="42.50"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  +4.02%  '
$ws.Range('D38').Formula = '="0.424"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +15.29%  '
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').Value = '0.0₃0777'
$ws.Range('E40').Value = '  +6.54%  '
$ws.Range('D41').Formula = '="0.139"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  +8.12%  '
$ws.Range('D42').Value = '3.242.47'
$ws.Range('E42').Value = '  +13.84%  '
$ws.Range('E43').Value = '  +14.43%  '
$ws.Range('D44').Formula = '="2.81"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +18.26%  '
$ws.Range('D45').Formula = '="0.999"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('E46').Value = '  +38.83%  '
$ws.Range('D47').Formula = '="0.0419"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  +4.37%  '
$ws.Range('E48').Value = '  +11.38%  '
$ws.Range('E49').Value = '  +9.86%  '
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('D51').Formula = '="3.08"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +4.76%  '
$excel.CutCopyMode = $false
